$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.012.68'
$ws.Range('E2').Value = '  -1.46%  '
$ws.Range('D3').Value = '1.782.31'
$ws.Range('E3').Value = '  -3.07%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '223.88'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.07%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.547'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.33%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '32.24'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.24%  '
$ws.Range('E9').Value = '  -3.46%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0702'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.95%  '
$ws.Range('E11').Value = '  -0.42%  '
$ws.Range('D12').Value = '2.038.06'
$ws.Range('E12').Value = '  -3.47%  '
$ws.Range('D13').Value = '1.789.54'
$ws.Range('E13').Value = '  -2.96%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.76'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.52%  '
$ws.Range('B15').Value = 'WrappedBTC'
$ws.Range('C15').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D15').Value = '34.001.02'
$ws.Range('E15').Value = '  -1.58%  '
$ws.Range('B16').Value = 'Polygon'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.620'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -4.49%  '
$ws.Range('E17').Value = '  -4.87%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '67.58'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -3.15%  '
$ws.Range('E19').Value = '  -4.09%  '
$ws.Range('E20').Value = '  -2.82%  '
$ws.Range('E21').Value = '  +0.19%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.63'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -5.85%  '
$ws.Range('E23').Value = '  -5.27%  '
$ws.Range('E24').Value = '  -2.47%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '159.31'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.55%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '16.25'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.34%  '
$ws.Range('E27').Value = '  -3.22%  '
$ws.Range('E28').Value = '  -2.75%  '
$ws.Range('E29').Value = '  +0.09%  '
$ws.Range('E30').Value = '  -4.62%  '
$ws.Range('E31').Value = '  -0.33%  '
$ws.Range('E32').Value = '  -4.36%  '
$ws.Range('E33').Value = '  -4.36%  '
$ws.Range('E34').Value = '  -7.22%  '
$ws.Range('D35').Value = '1.390.08'
$ws.Range('E35').Value = '  -4.64%  '
$ws.Range('E36').Value = '  -2.29%  '
$ws.Range('E37').Value = '  -3.34%  '
$ws.Range('E38').Value = '  -4.58%  '
$ws.Range('E39').Value = '  -0.90%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.19'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.23%  '
$ws.Range('E41').Value = '  -3.67%  '
$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '78.17'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -5.74%  '
$ws.Range('B43').Value = 'ARBITRUM'
$ws.Range('C43').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.908'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -7.13%  '
$ws.Range('D44').Value = '0.0₆0146'
$ws.Range('E44').Value = '  +13.64%  '
$ws.Range('E45').Value = '  +1.26%  '
$ws.Range('E46').Value = '  -0.14%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '106.77'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.35%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.84'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -5.15%  '
$ws.Range('E49').Value = '  -0.35%  '
$ws.Range('D50').Value = '1.938.37'
$ws.Range('E50').Value = '  -3.66%  '
$ws.Range('E51').Value = '  -0.34%  '
